$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply each updated cell as a literal text value (apostrophe-prefixed
# to stop Excel re-interpreting dotted/numeric-looking strings as numbers),
# then reset the style back to Normal so no stray NumberFormat/quote-prefix
# style gets attached to the cell (matches original "no style" cells).
function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.Value = "'" + $val
    $c.Style = "Normal"
}

Set-TextValue "D2" "27.851.10"
Set-TextValue "E2" "  +2.50%  "
Set-TextValue "D3" "1.871.94"
Set-TextValue "E3" "  +0.80%  "
Set-TextValue "D4" "1.012"
Set-TextValue "E4" "  -0.63%  "
Set-TextValue "D5" "313.36"
Set-TextValue "E5" "  +0.74%  "
Set-TextValue "D6" "1.011"
Set-TextValue "E6" "  -0.63%  "
Set-TextValue "D7" "0.4828"
Set-TextValue "E7" "  +0.79%  "
Set-TextValue "D8" "0.3817"
Set-TextValue "E8" "  +2.93%  "
Set-TextValue "D9" "0.07376"
Set-TextValue "D10" "0.9401"
Set-TextValue "E10" "  +0.43%  "
Set-TextValue "E11" "  +4.97%  "
Set-TextValue "D12" "0.07798"
Set-TextValue "E12" "  -0.77%  "
Set-TextValue "D13" "1.863.38"
Set-TextValue "E13" "  -0.25%  "
Set-TextValue "D14" "5.493"
Set-TextValue "E14" "  +1.43%  "
Set-TextValue "D15" "6.610"
Set-TextValue "E15" "  +1.23%  "
Set-TextValue "D16" "91.17"
Set-TextValue "E16" "  +1.44%  "
Set-TextValue "D17" "1.013"
Set-TextValue "E17" "  -0.61%  "
Set-TextValue "D18" "0.000008863"
Set-TextValue "E18" "  +1.65%  "
Set-TextValue "E19" "  -0.54%  "
Set-TextValue "D20" "27.858.42"
Set-TextValue "E20" "  +2.34%  "
Set-TextValue "E21" "  +1.00%  "
Set-TextValue "D22" "5.117"
Set-TextValue "E22" "  +0.45%  "
Set-TextValue "D23" "2.118.19"
Set-TextValue "E23" "  +1.68%  "
Set-TextValue "D24" "10.86"
Set-TextValue "E24" "  +1.73%  "
Set-TextValue "D25" "1.947"
Set-TextValue "E25" "  +0.31%  "
Set-TextValue "D26" "157.12"
Set-TextValue "E26" "  +2.36%  "
Set-TextValue "E27" "  +0.49%  "
Set-TextValue "E28" "  +2.91%  "
Set-TextValue "E29" "  +0.43%  "
Set-TextValue "D30" "4.984"
Set-TextValue "E30" "  +1.11%  "
Set-TextValue "D31" "0.08900"
Set-TextValue "E31" "  +0.18%  "
Set-TextValue "D32" "3.334"
Set-TextValue "E32" "  +0.68%  "
Set-TextValue "D33" "1.228"
Set-TextValue "E33" "  +3.65%  "
Set-TextValue "D34" "0.7702"
Set-TextValue "E34" "  +4.41%  "
Set-TextValue "D35" "4.664"
Set-TextValue "E35" "  +2.15%  "
Set-TextValue "D36" "2.747"
Set-TextValue "E36" "  +2.64%  "
Set-TextValue "D37" "1.132"
Set-TextValue "E37" "  +1.23%  "
Set-TextValue "E38" "  +2.02%  "
Set-TextValue "D39" "0.5606"
Set-TextValue "E39" "  +5.21%  "
Set-TextValue "E40" "  +1.91%  "
Set-TextValue "E41" "  +0.26%  "
Set-TextValue "D42" "7.051"
Set-TextValue "D43" "8.552"
Set-TextValue "E43" "  +2.49%  "
Set-TextValue "D44" "0.1528"
Set-TextValue "E44" "  +0.11%  "
Set-TextValue "D45" "0.4876"
Set-TextValue "E45" "  +2.27%  "
Set-TextValue "D46" "10.68"
Set-TextValue "E46" "  +0.40%  "
Set-TextValue "B47" "Quant"
Set-TextValue "C47" "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue "D47" "105.37"
Set-TextValue "E47" "  +2.87%  "
Set-TextValue "B48" "PaxDollar"
Set-TextValue "C48" "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextValue "D48" "1.012"
Set-TextValue "E48" "  -0.61%  "
Set-TextValue "E49" "  +2.04%  "
Set-TextValue "D50" "68.08"
Set-TextValue "E50" "  +2.64%  "
Set-TextValue "D51" "0.06126"
Set-TextValue "E51" "  +0.88%  "
